# Fix bugs: - Read formula
# Rename "Sheet2" to "Data", populate it with sample data covering
# int/double/string/percent/date values plus formulas (including a
# volatile TODAY() read-formula scenario), and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$ws2.Name = "Data"
$ws2.Activate()

# Row 1 - Int
$ws2.Range("A1").Value = "Int"
$ws2.Range("B1").Value = 1

# Row 2 - Double
$ws2.Range("A2").Value = "Double"
$ws2.Range("B2").Value = 0.25

# Row 3 - Double Formula
$ws2.Range("A3").Value = "Double Formula"
$ws2.Range("B3").Formula = "=6/10"

# Row 4 - String
$ws2.Range("A4").Value = "String"
$ws2.Range("B4").Value = "Test"

# Row 5 - Percent
$ws2.Range("A5").Value = "Percent"
$ws2.Range("B5").NumberFormat = "0%"
$ws2.Range("B5").Value = 0.1

# Row 6 - String Formular
$ws2.Range("A6").Value = "String Formular"
$ws2.Range("B6").Formula = '="A" & "B"'

# Row 7 - Date
$ws2.Range("A7").Value = "Date"
$ws2.Range("B7").NumberFormat = "d-mmm"
$ws2.Range("B7").Value = 40939

# Labels for rows 8/9 - set "Date Formular" before "Date Func" so the
# shared-string table order matches the workbook that was read back in.
$ws2.Range("A9").Value = "Date Formular"
$ws2.Range("A8").Value = "Date Func"

# Row 8 - Date Func (volatile TODAY() read formula)
# Row 9 - Date Formular (reads row 8 + 1)
# Apply the number format first and copy it across so both cells share
# a single style record, then fill in the formulas.
$ws2.Range("B8").NumberFormat = "mm-dd-yy"
$ws2.Range("B8").Copy() | Out-Null
$ws2.Range("B9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws2.Range("B8").Formula = "=TODAY()"
$ws2.Range("B9").Formula = "=B8+1"

# Column widths (closest achievable values to the authored widths).
$ws2.Columns.Item(1).ColumnWidth = 14.5
$ws2.Columns.Item(2).ColumnWidth = 13

# Restore the cursor/selection position recorded on the sheet.
$ws2.Range("G19").Select() | Out-Null
